$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Insert a new slide at position 10 ("Title and Content" layout), delete
#    its (unused) Title placeholder, and fill the Content placeholder with
#    the new "fig 5 runs" notes.  This pushes the three trailing slides
#    (old positions 10-12) down to 11-13, matching the new p:sldIdLst order.
# ---------------------------------------------------------------------------
$s = $p.Slides.Add(10, 2)
$s.Shapes.Item(1).Delete()
$content = $s.Shapes.Item(1)
$tr = $content.TextFrame.TextRange
$dash = [char]0x2013

$tr.Text = "Main source of variation is initial proportion"
$tr.InsertAfter("`rThree parameters: initial proportion, number of cell types, difference in cell types") | Out-Null
$tr.InsertAfter("`rA $dash show variance with one type is small") | Out-Null
$tr.InsertAfter("`rB $dash show variance is large with two cell types and varying prop") | Out-Null
$tr.InsertAfter("`rC $dash show adding cell types does not increase this variance") | Out-Null
$tr.InsertAfter("`rD $dash show increasing variance with increase variance in prop") | Out-Null
$tr.InsertAfter("`rEach fig has 3 components: total over time, prop ") | Out-Null
$tr.InsertAfter("dist") | Out-Null
$tr.InsertAfter(", total ") | Out-Null
$tr.InsertAfter("dist") | Out-Null
$tr.InsertAfter(" at 48 hours") | Out-Null
$tr.InsertAfter("`r") | Out-Null

# ---------------------------------------------------------------------------
# 2. Bump the cached "datetimeFigureOut" date field from 10/10/2017 to
#    10/21/2017 everywhere it is cached: the slide master and every slide
#    layout's Date Placeholder.
# ---------------------------------------------------------------------------
$newDate = "10/21/2017"
$m = $p.SlideMaster

for ($j = 1; $j -le $m.Shapes.Count; $j++) {
    $sh = $m.Shapes.Item($j)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($i = 1; $i -le $m.CustomLayouts.Count; $i++) {
    $cl = $m.CustomLayouts.Item($i)
    for ($j = 1; $j -le $cl.Shapes.Count; $j++) {
        $sh = $cl.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

Write-Output "done"
